$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Heure" -> "Heure " (trailing space)
$ws.Range("C2").Value = "Heure "

# Row 11 data moves up to row 8 (Côme / Mathieu), dropping the time column
$ws.Range("B8").Value = "Côme"
$ws.Range("D8").Value = "Mathieu"

# Row 11 data (second half) moves up to row 9 (Côme / Didier)
$ws.Range("B9").Value = "Côme"
$ws.Range("D9").Value = "Didier"

# Row 12 data moves up to row 10 (Côme / Baptiste)
$ws.Range("B10").Value = "Côme"
$ws.Range("D10").Value = "Baptiste"

# Clear out what used to be rows 11 and 12 entirely (they no longer exist)
$ws.Range("A11:E12").Clear()
